# The presentation currently has its theme parts "swapped" relative to the
# target state: ppt/theme/theme2.xml (the theme actually driving the slide
# master / every slide) holds the "Integral" / "Red Violet" colour scheme,
# while ppt/theme/theme1.xml (only loosely linked from the notes master)
# holds the "Office Theme" colour scheme. The authored edit swaps the two
# parts' contents so that the master's theme becomes "Office Theme" colours
# (and the other part becomes "Integral").
#
# The live, editable theme reachable from the PowerPoint object model is the
# one attached to the (single) slide master/design - that is the part that
# actually renders on every slide. We drive the swap by writing the twelve
# "Office Theme" scheme colours onto that theme's ThemeColorScheme, in the
# standard PowerPoint order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
